# New crime data collected: update weekly CompStat figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: bump the Volume/Number and the reporting week dates.
# Use Characters(...).Text so the existing rich-text run (font/size/color)
# is preserved and only the digits inside the run are replaced.
# ---------------------------------------------------------------------------

# A8: "Volume 30   Number  39" -> "Volume 30   Number  40"
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "40"

# C9: "Report Covering the Week  9/25/2023  Through  10/1/2023"
#  -> "Report Covering the Week  10/2/2023  Through  10/8/2023"
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "10/2/2023"
$c9.Characters(47, 8).Text = "10/8/2023"

# ---------------------------------------------------------------------------
# Helper: convert a numeric cell to the "N/A" text marker cell (style s=14),
# reusing an untouched same-column donor cell (row 14) purely for its format.
# ---------------------------------------------------------------------------
function Set-TextMarker($target, $donorAddr, $text) {
    $t = $ws.Range($target)
    $t.NumberFormat = "@"
    $t.Value = $text
    $ws.Range($donorAddr).Copy() | Out-Null
    $t.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Helper: convert a text marker cell back into a plain numeric cell, reusing
# an untouched same-column donor cell purely for its (numeric) format.
function Set-NumericValue($target, $donorAddr, $value) {
    $t = $ws.Range($target)
    $t.Value = $value
    $ws.Range($donorAddr).Copy() | Out-Null
    $t.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -53.333333333333

# ---------------------------------------------------------------------------
# Row 16 (Robbery) - D16/E16 become text "N/A" markers (s=14)
# ---------------------------------------------------------------------------
Set-TextMarker "D16" "D14" "0"
Set-TextMarker "E16" "E14" "***.*"
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 113
$ws.Range("K16").Value = 39.506172839506
$ws.Range("L16").Value = 109.259259259259
$ws.Range("M16").Value = -10.317460317460
$ws.Range("N16").Value = -81.166666666666

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 9.090909090909
$ws.Range("I17").Value = 139
$ws.Range("J17").Value = 121
$ws.Range("K17").Value = 14.876033057851
$ws.Range("L17").Value = 23.008849557522
$ws.Range("M17").Value = 52.747252747252
$ws.Range("N17").Value = -37.668161434977

# ---------------------------------------------------------------------------
# Row 18 (Burglary) - C18 becomes a plain number again (s=15)
# ---------------------------------------------------------------------------
Set-NumericValue "C18" "C19" 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -64.285714285714
$ws.Range("I18").Value = 73
$ws.Range("J18").Value = 98
$ws.Range("K18").Value = -25.510204081632
$ws.Range("L18").Value = 69.767441860465
$ws.Range("M18").Value = 23.728813559322
$ws.Range("N18").Value = -84.434968017057

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 52.380952380952
$ws.Range("I19").Value = 310
$ws.Range("J19").Value = 267
$ws.Range("K19").Value = 16.104868913857
$ws.Range("L19").Value = 64.021164021164
$ws.Range("M19").Value = 49.758454106280
$ws.Range("N19").Value = -40.384615384615

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 42.857142857142
$ws.Range("I20").Value = 86
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = 48.275862068965
$ws.Range("L20").Value = 126.315789473684
$ws.Range("M20").Value = 330
$ws.Range("N20").Value = -70.748299319727

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 17
$ws.Range("E21").Value = 6.25
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 11.290322580645
$ws.Range("I21").Value = 728
$ws.Range("J21").Value = 641
$ws.Range("K21").Value = 13.572542901716
$ws.Range("L21").Value = 64.705882352941
$ws.Range("M21").Value = 40.812379110251
$ws.Range("N21").Value = -65.933551708001

# ---------------------------------------------------------------------------
# Row 22 (Transit) - C22 becomes text "0" marker, D22/E22 become plain numbers
# ---------------------------------------------------------------------------
Set-TextMarker "C22" "C14" "0"
Set-NumericValue "D22" "D19" 1
Set-NumericValue "E22" "E19" -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = 93.333333333333
$ws.Range("L22").Value = 26.086956521739

# ---------------------------------------------------------------------------
# Row 23 (Housing)
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 116
$ws.Range("J23").Value = 91
$ws.Range("K23").Value = 27.472527472527
$ws.Range("L23").Value = 36.470588235294
$ws.Range("M23").Value = 87.096774193548

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = -12.765957446808
$ws.Range("I24").Value = 396
$ws.Range("J24").Value = 401
$ws.Range("K24").Value = -1.246882793017
$ws.Range("L24").Value = 27.331189710610
$ws.Range("M24").Value = -18.852459016393

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 186
$ws.Range("J25").Value = 165
$ws.Range("K25").Value = 12.727272727272
$ws.Range("L25").Value = 17.721518987341
$ws.Range("M25").Value = -19.480519480519

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*) - C26 becomes a plain number again (s=15)
# ---------------------------------------------------------------------------
Set-NumericValue "C26" "C19" 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = -54.545454545454
$ws.Range("L26").Value = 42.857142857142

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 29.032258064516

# ---------------------------------------------------------------------------
# Row 30 (Hate Crimes) - C30/G30 become text "0" markers, H30 becomes "***.*"
# ---------------------------------------------------------------------------
Set-TextMarker "C30" "C14" "0"
Set-TextMarker "G30" "G14" "0"
Set-TextMarker "H30" "H14" "***.*"

Write-Output "done"
